$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 (B23:C23) currently holds the "LOT2028 ... (Requisito fraco)" text,
# row 24 (B24:C24) currently holds the "LOT2052 ... (Indicação de Conjunto)" text.
# The shared-string table order needs to flip so the LOT2052 entry precedes the
# LOT2028 entry; since the sheet keeps referencing the same cells (B23:C23 then
# B24:C24), the net visible effect is that the two rows' text values swap.

$textLOT2052 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"
$textLOT2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

$ws.Range("B23").Value = $textLOT2052
$ws.Range("C23").Value = $textLOT2052

$ws.Range("B24").Value = $textLOT2028
$ws.Range("C24").Value = $textLOT2028
